$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Text)
    $helper = $ws.Range("Z1000")
    $helper.NumberFormat = "@"
    $helper.Value = $Text
    $helper.Copy() | Out-Null
    $ws.Range($CellRef).PasteSpecial(-4163) | Out-Null
    $helper.ClearContents()
    $ws.Application.CutCopyMode = $false
}

$ws.Range('D2').Value = '71.520.81'
$ws.Range('E2').Value = '  +2.14%  '
$ws.Range('D3').Value = '3.814.45'
$ws.Range('E3').Value = '  +0.42%  '
Set-TextValue -CellRef 'D4' -Text '0.999'
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue -CellRef 'D5' -Text '702.69'
$ws.Range('E5').Value = '  +5.95%  '
Set-TextValue -CellRef 'D6' -Text '174.52'
$ws.Range('E6').Value = '  +4.77%  '
$ws.Range('D7').Value = '3.815.90'
$ws.Range('E7').Value = '  +0.49%  '
Set-TextValue -CellRef 'D8' -Text '0.999'
$ws.Range('E8').Value = '  -0.07%  '
Set-TextValue -CellRef 'D9' -Text '0.530'
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('E10').Value = '  +2.44%  '
Set-TextValue -CellRef 'D11' -Text '7.28'
$ws.Range('E11').Value = '  +4.34%  '
Set-TextValue -CellRef 'D12' -Text '0.463'
$ws.Range('E12').Value = '  +0.71%  '
Set-TextValue -CellRef 'D13' -Text '0.0000259'
$ws.Range('E13').Value = '  +7.07%  '
Set-TextValue -CellRef 'D14' -Text '36.41'
$ws.Range('E14').Value = '  +2.02%  '
$ws.Range('D15').Value = '4.451.79'
$ws.Range('E15').Value = '  +0.33%  '
$ws.Range('D16').Value = '3.805.18'
$ws.Range('D17').Value = '71.378.53'
$ws.Range('E17').Value = '  +1.98%  '
Set-TextValue -CellRef 'D18' -Text '17.73'
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('E21').Value = '  +6.99%  '
Set-TextValue -CellRef 'D22' -Text '483.20'
Set-TextValue -CellRef 'D23' -Text '0.715'
$ws.Range('E23').Value = '  +0.21%  '
Set-TextValue -CellRef 'D24' -Text '84.64'
$ws.Range('E24').Value = '  +2.31%  '
$ws.Range('E25').Value = '  -0.68%  '
Set-TextValue -CellRef 'D26' -Text '12.37'
$ws.Range('E26').Value = '  +1.13%  '
Set-TextValue -CellRef 'D27' -Text '10.62'
$ws.Range('E27').Value = '  +2.91%  '
$ws.Range('E28').Value = '  +1.96%  '
$ws.Range('D29').Value = '3.964.98'
$ws.Range('E29').Value = '  +0.30%  '
Set-TextValue -CellRef 'D30' -Text '3.16'
$ws.Range('E30').Value = '  +12.40%  '
$ws.Range('E31').Value = '  +0.05%  '
Set-TextValue -CellRef 'D32' -Text '7.64'
$ws.Range('E32').Value = '  +3.81%  '
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('E34').Value = '  +5.39%  '
Set-TextValue -CellRef 'D35' -Text '29.59'
$ws.Range('E35').Value = '  +1.39%  '
Set-TextValue -CellRef 'D36' -Text '9.30'
$ws.Range('E36').Value = '  +2.71%  '
Set-TextValue -CellRef 'D37' -Text '1.00'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('E38').Value = '  +1.75%  '
Set-TextValue -CellRef 'D39' -Text '3.48'
$ws.Range('E39').Value = '  +5.52%  '
Set-TextValue -CellRef 'D40' -Text '6.04'
$ws.Range('E40').Value = '  +2.31%  '
Set-TextValue -CellRef 'D41' -Text '2.31'
$ws.Range('E41').Value = '  +12.61%  '
Set-TextValue -CellRef 'D42' -Text '0.990'
$ws.Range('E42').Value = '  +2.48%  '
Set-TextValue -CellRef 'D43' -Text '0.997'
$ws.Range('E43').Value = '  -0.35%  '
Set-TextValue -CellRef 'D45' -Text '0.000319'
$ws.Range('E45').Value = '  +16.92%  '
Set-TextValue -CellRef 'D46' -Text '165.13'
$ws.Range('E46').Value = '  +3.87%  '
Set-TextValue -CellRef 'D47' -Text '44.94'
$ws.Range('E47').Value = '  -0.23%  '
Set-TextValue -CellRef 'D48' -Text '48.77'
$ws.Range('E48').Value = '  +1.58%  '
Set-TextValue -CellRef 'D49' -Text '0.304'
$ws.Range('E49').Value = '  +1.15%  '
Set-TextValue -CellRef 'D50' -Text '416.95'
$ws.Range('E50').Value = '  +7.32%  '
$ws.Range('E51').Value = '  -1.53%  '
